# Edit 1 (slide 2, "STGraph - Overview"):
# Split the combined bullet "Different layouts ... -> Yet many modern ..."
# into two separate paragraphs: the first stays at the existing indent
# level (0) and keeps the text up to the ellipsis; the second is a new,
# more-indented (level 1) paragraph holding the remainder of the sentence.
$p = $ppt.ActivePresentation

$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange

# Touch the paragraph's text first with an unrelated placeholder so the
# later assignment of the real (shorter, prefix-matching) text is not
# treated as a partial in-place edit that fragments the run in two.
$tr2.Paragraphs(9).Text = "placeholder"
$tr2.Paragraphs(9).Text = "Different layouts ➔ different ingestion and workload capabilities…"

# Insert a brand-new paragraph right after it containing the second half
# of the original sentence, then demote it to indent level 1 (lvl="1").
[void]$tr2.Paragraphs(9).InsertAfter([char]13 + "Yet many modern data-intensive applications combine both, e.g., IoT systems, Digital Twins, and pervasive computing.")
$tr2.Paragraphs(10).IndentLevel = 2

# Edit 2 (slide 6, "STGraph - Operations"):
# Reword the "Naive nested-Loop join strategy;" bullet.
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange

$tr6.Paragraphs(12).Text = "placeholder"
$tr6.Paragraphs(12).Text = "Support for join operations through naive nested-Loop join;"
